# "Change in project name" - the "yes" execute-flag is flipped to "no"
# for a couple of rows, and the previously active RUNMANAGER sheet/cell
# selection is swapped for the DATA sheet/cell selection.

$wb = $excel.ActiveWorkbook

$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$data       = $wb.Worksheets.Item("DATA")

# RUNMANAGER!C2 : "yes" -> "no"
$runmanager.Range("C2").Value = "no"

# DATA!B4 and DATA!B5 : "yes" -> "no"
$data.Range("B4").Value = "no"
$data.Range("B5").Value = "no"

# Update the remembered selection on each sheet ...
$runmanager.Range("B10").Select()
$data.Range("C13").Select()

# ... and make DATA the active (visible) tab, as it is last selected.
$data.Activate()
